$wb = $excel.ActiveWorkbook

# ---------- Sheet: 展览 ----------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 840
$ws1.Range("F9").Value = 125
$ws1.Range("F11").Value = 1229
$ws1.Range("F15").Value = 897
$ws1.Range("F18").Value = 81
$ws1.Range("F20").Value = 817
$ws1.Range("F21").Value = 1757
$ws1.Range("F22").Value = 3172
$ws1.Range("F23").Value = 931
$ws1.Range("F25").Value = 2320
$ws1.Range("F27").Value = 12
$ws1.Range("F28").Value = 3188
$ws1.Range("F29").Value = 659
$ws1.Range("F30").Value = 810
$ws1.Range("F31").Value = 22
$ws1.Range("F32").Value = 95
$ws1.Range("F33").Value = 746
$ws1.Range("F35").Value = 143
$ws1.Range("F36").Value = 75
$ws1.Range("F38").Value = 1139
$ws1.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202404/43sjLXZh1712910203022.jpeg"
$ws1.Range("F39").Value = 1820
$ws1.Range("F40").Value = 419
$ws1.Range("F43").Value = 214
$ws1.Range("F44").Value = 140
$ws1.Range("F45").Value = 194
$ws1.Range("F46").Value = 57

# ---------- Sheet: 演出 ----------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 144
$ws2.Range("F12").Value = 96

# ---------- Sheet: 本地生活 ----------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 116

# ---------- Sheet: 全部类型 ----------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 840
$ws4.Range("F7").Value = 125
$ws4.Range("F8").Value = 1229
$ws4.Range("F11").Value = 897
$ws4.Range("F12").Value = 144
$ws4.Range("F16").Value = 81
$ws4.Range("F17").Value = 817
$ws4.Range("F18").Value = 1757
$ws4.Range("F19").Value = 3172
$ws4.Range("F20").Value = 931
$ws4.Range("F23").Value = 2320
$ws4.Range("F24").Value = 12
$ws4.Range("F25").Value = 3188
$ws4.Range("F26").Value = 659
$ws4.Range("F27").Value = 810
$ws4.Range("F29").Value = 22
$ws4.Range("F33").Value = 95
$ws4.Range("F34").Value = 96
$ws4.Range("F35").Value = 746
$ws4.Range("F37").Value = 143
$ws4.Range("F38").Value = 75
$ws4.Range("F41").Value = 1139
$ws4.Range("I41").Value = "//i0.hdslb.com/bfs/openplatform/202404/43sjLXZh1712910203022.jpeg"
$ws4.Range("F42").Value = 1820
$ws4.Range("F44").Value = 419
$ws4.Range("F46").Value = 214
$ws4.Range("F47").Value = 140
$ws4.Range("F48").Value = 194
$ws4.Range("F49").Value = 57
